$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was inserted into the weekly series at row 100,
# pushing every following record (old rows 100-169) down by one row
# (new rows 101-170). Dimension grows from A1:R169 to A1:R170.
$ws.Rows.Item(100).Insert()

$ws.Cells.Item(100, 1).Value = 5
$ws.Cells.Item(100, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(100, 3).Value = "Maule"
$ws.Cells.Item(100, 4).Value = 44574
$ws.Cells.Item(100, 5).Value = 7
$ws.Cells.Item(100, 6).Value = 100112024
$ws.Cells.Item(100, 7).Value = "Choclo"
$ws.Cells.Item(100, 8).Value = "Choclero"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 50000
$ws.Cells.Item(100, 11).Value = 230
$ws.Cells.Item(100, 12).Value = 250
$ws.Cells.Item(100, 13).Value = 238
$ws.Cells.Item(100, 14).Value = "`$/unidad"
$ws.Cells.Item(100, 15).Value = "Región del Maule"
$ws.Cells.Item(100, 16).Value = 238
$ws.Cells.Item(100, 17).Value = 1
$ws.Cells.Item(100, 18).Value = "Hortaliza"
